$wb = $excel.ActiveWorkbook

$wsMass    = $wb.Worksheets.Item("Rough Mass Calculations")
$wsFlight  = $wb.Worksheets.Item("Rough Flight Calcs")

# --- Flight Sim calcs for SAC 2018 ---
# The "Rough Flight Calcs" sheet pulls its launch-mass figure (column D,
# rows 3:5) from the "Rough Mass Calculations" sheet. Previously all three
# scenarios referenced the same absolute cell ($D$26). Update them to pull
# from the three different rocket configurations (D26/D27/D28) and subtract
# 2 kg (e.g. motor/propellant mass already accounted for elsewhere).
$wsFlight.Range("D3").Formula = "=('Rough Mass Calculations'!D26/1000) - 2"
$wsFlight.Range("D4").Formula = "=('Rough Mass Calculations'!D27/1000) - 2"
$wsFlight.Range("D5").Formula = "=('Rough Mass Calculations'!D28/1000) - 2"

# --- View/selection bookkeeping to mirror the author's saved workbook state ---
# Scroll "Rough Mass Calculations" back to the left edge (A1) instead of U1.
$wsMass.Activate() | Out-Null
$wsMass.Range("A1").Select() | Out-Null
$wsMass.Range("AG6").Select() | Out-Null

# "Rough Flight Calcs" ends up the active sheet, zoomed in further, scrolled
# right to column R, with U9 selected.
$wsFlight.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 18
$excel.ActiveWindow.ScrollRow = 1
$wsFlight.Range("U9").Select() | Out-Null
$excel.ActiveWindow.Zoom = 136
